# Update NATMI TPM-derived values for the Ccl25-Ccr9 LR-pair sheet.
#
# Columns (1-indexed):
#   A Sending cluster           G Ligand avg expr value      M Receptor avg expr value
#   B Ligand symbol              H Ligand total expr value    N Receptor total expr value
#   C Receptor symbol            I Ligand specificity (avg)   O Receptor specificity (avg)
#   D Target cluster             J Ligand specificity (total) P Receptor specificity (total)
#                                                               Q Edge avg expr weight  (=G*M)
#                                                               R Edge total expr weight (=H*N)
#                                                               S Edge avg expr specificity  (=I*O)
#                                                               T Edge total expr specificity (=J*P)
#
# G/H only depend on the Sending cluster (col A); M/N only depend on the
# Target cluster (col D). I/J/O/P are each value normalised against the sum
# of that same column across all (five) clusters. Q/R/S/T are simple
# products of the ligand/receptor figures for that row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New per-cluster ligand expression values (average, total) -- column G, H
$GH = @{}
$GH["ECs"]               = @(3.505502,            10.516506)
$GH["FAPs"]              = @(6.985602333333333,   20.956807)
$GH["Inflammatory-Mac"]  = @(4.346963,            13.040889)
$GH["MuSCs"]             = @(6.487399,            12.974798)
$GH["Resolving-Mac"]     = @(5.098167333333333,   15.294502)

# New per-cluster receptor expression values (average, total) -- column M, N
$MN = @{}
$MN["ECs"]               = @(0.4829603333333334,  1.448881)
$MN["FAPs"]              = @(1.419989,            4.259967)
$MN["Inflammatory-Mac"]  = @(4.546560333333333,   13.639681)
$MN["MuSCs"]             = @(0.8193049999999999,  1.63861)
$MN["Resolving-Mac"]     = @(1.758262666666667,   5.274788000000001)

$firstRow = 2
$lastRow = 26

# Write the new ligand (G,H) / receptor (M,N) figures for every data row,
# keyed off that row's Sending cluster (A) / Target cluster (D).
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $sending = $ws.Cells.Item($r, 1).Value2
    $target  = $ws.Cells.Item($r, 4).Value2

    $g = $GH[$sending][0]
    $h = $GH[$sending][1]
    $m = $MN[$target][0]
    $n = $MN[$target][1]

    $ws.Cells.Item($r, 7).Value  = $g    # G
    $ws.Cells.Item($r, 8).Value  = $h    # H
    $ws.Cells.Item($r, 13).Value = $m    # M
    $ws.Cells.Item($r, 14).Value = $n    # N
}

# Ligand/receptor derived-specificity columns are each figure divided by the
# sum of that figure across all five clusters -- recompute the sums from the
# values we just wrote so every row sees the same (new) totals.
$sumG = 0.0
$sumH = 0.0
foreach ($k in $GH.Keys) {
    $sumG = $sumG + $GH[$k][0]
    $sumH = $sumH + $GH[$k][1]
}

$sumM = 0.0
$sumN = 0.0
foreach ($k in $MN.Keys) {
    $sumM = $sumM + $MN[$k][0]
    $sumN = $sumN + $MN[$k][1]
}

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $sending = $ws.Cells.Item($r, 1).Value2
    $target  = $ws.Cells.Item($r, 4).Value2

    $g = $GH[$sending][0]
    $h = $GH[$sending][1]
    $m = $MN[$target][0]
    $n = $MN[$target][1]

    $i = $g / $sumG     # Ligand derived specificity of average expression value
    $j = $h / $sumH     # Ligand derived specificity of total expression value
    $o = $m / $sumM     # Receptor derived specificity of average expression value
    $p = $n / $sumN     # Receptor derived specificity of total expression value

    $ws.Cells.Item($r, 9).Value  = $i    # I
    $ws.Cells.Item($r, 10).Value = $j    # J
    $ws.Cells.Item($r, 15).Value = $o    # O
    $ws.Cells.Item($r, 16).Value = $p    # P

    $q = $g * $m        # Edge average expression weight
    $rr = $h * $n       # Edge total expression weight
    $s = $i * $o        # Edge average expression derived specificity
    $t = $j * $p        # Edge total expression derived specificity

    $ws.Cells.Item($r, 17).Value = $q    # Q
    $ws.Cells.Item($r, 18).Value = $rr   # R
    $ws.Cells.Item($r, 19).Value = $s    # S
    $ws.Cells.Item($r, 20).Value = $t    # T
}

Write-Output "updated rows $firstRow..$lastRow"
